$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new EUR->ARS rate reading as row 29.
$row = 29

# Column A holds a date-like string (e.g. "2025-09-20"). Force the cell to
# Text format before assigning so Excel stores it as plain text rather than
# auto-converting it into a date serial value, matching the existing rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-20"
# Reset the cell style back to the default ("Normal") so no extra style index
# lingers on the cell, keeping it consistent with the rest of the sheet.
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "15:18:22"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,777.8410"
